$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

$ws.Range("D2").Value = "29.847.82"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").Value = "1.869.57"
$ws.Range("E3").Value = "  -1.38%  "

Set-TextValue $ws "D4" "0.9998"
$ws.Range("E4").Value = "  -0.10%  "

Set-TextValue $ws "D5" "0.7381"
$ws.Range("E5").Value = "  -4.74%  "

Set-TextValue $ws "D6" "241.88"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("E7").Value = "  -0.11%  "

Set-TextValue $ws "D8" "0.3149"
$ws.Range("E8").Value = "  +0.26%  "

Set-TextValue $ws "D9" "24.65"
$ws.Range("E9").Value = "  -4.52%  "

Set-TextValue $ws "D10" "0.07099"
$ws.Range("E10").Value = "  -2.10%  "

Set-TextValue $ws "D11" "0.08375"
$ws.Range("E11").Value = "  -5.68%  "

Set-TextValue $ws "D12" "0.7518"
$ws.Range("E12").Value = "  -2.98%  "

Set-TextValue $ws "D13" "5.443"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "1.866.42"
$ws.Range("E14").Value = "  +0.73%  "

Set-TextValue $ws "D15" "92.44"
$ws.Range("E15").Value = "  -2.57%  "

$ws.Range("D16").Value = "29.864.29"
$ws.Range("E16").Value = "  +0.07%  "

Set-TextValue $ws "D17" "6.026"
$ws.Range("E17").Value = "  -2.86%  "

Set-TextValue $ws "D18" "13.57"
$ws.Range("E18").Value = "  -3.07%  "

Set-TextValue $ws "D19" "242.83"
$ws.Range("E19").Value = "  -1.61%  "

Set-TextValue $ws "D20" "0.000007816"
$ws.Range("E20").Value = "  -1.24%  "

Set-TextValue $ws "D21" "0.9986"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "2.116.29"
$ws.Range("E22").Value = "  -0.54%  "

Set-TextValue $ws "D23" "7.909"
$ws.Range("E23").Value = "  -3.13%  "

Set-TextValue $ws "D24" "0.9999"
$ws.Range("E24").Value = "  -0.11%  "

Set-TextValue $ws "D25" "0.1561"
$ws.Range("E25").Value = "  -1.73%  "

Set-TextValue $ws "D26" "9.296"
$ws.Range("E26").Value = "  -2.78%  "

Set-TextValue $ws "D27" "164.19"
$ws.Range("E27").Value = "  +0.74%  "

Set-TextValue $ws "D28" "18.57"
$ws.Range("E28").Value = "  -1.55%  "

Set-TextValue $ws "D29" "2.016"
$ws.Range("E29").Value = "  -1.78%  "

Set-TextValue $ws "D30" "1.474"
$ws.Range("E30").Value = "  +3.28%  "

Set-TextValue $ws "D31" "4.658"
$ws.Range("E31").Value = "  +2.84%  "

$ws.Range("E32").Value = "  -1.37%  "

Set-TextValue $ws "D33" "4.298"
$ws.Range("E33").Value = "  +4.28%  "

Set-TextValue $ws "D34" "0.05316"
$ws.Range("E34").Value = "  -3.91%  "

Set-TextValue $ws "D35" "1.233"
$ws.Range("E35").Value = "  -1.30%  "

Set-TextValue $ws "D36" "0.7533"
$ws.Range("E36").Value = "  -0.31%  "

Set-TextValue $ws "D37" "1.001"
$ws.Range("E37").Value = "  +0.08%  "

Set-TextValue $ws "D38" "2.699"
$ws.Range("E38").Value = "  -0.93%  "

Set-TextValue $ws "D39" "0.01955"
$ws.Range("E39").Value = "  -0.74%  "

Set-TextValue $ws "D40" "2.750"
$ws.Range("E40").Value = "  -1.55%  "

$ws.Range("E41").Value = "  -1.11%  "

$ws.Range("D42").Value = "1.099.36"
$ws.Range("E42").Value = "  +0.80%  "

Set-TextValue $ws "D43" "6.071"
$ws.Range("E43").Value = "  -0.23%  "

Set-TextValue $ws "D44" "72.12"
$ws.Range("E44").Value = "  -2.71%  "

Set-TextValue $ws "D45" "0.8613"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("E46").Value = "  -0.02%  "

Set-TextValue $ws "D47" "103.11"
$ws.Range("E47").Value = "  +0.23%  "

Set-TextValue $ws "D48" "7.680"
$ws.Range("E48").Value = "  +0.68%  "

Set-TextValue $ws "D49" "1.841"
$ws.Range("E49").Value = "  -2.95%  "

Set-TextValue $ws "D50" "3.053"
$ws.Range("E50").Value = "  +1.94%  "

$ws.Range("D51").Value = "2.016.87"
$ws.Range("E51").Value = "  -1.67%  "
